$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
# Leading apostrophe forces these numeric-looking strings to stay text
# (matches source which stores them as inline strings, not numbers).
$ws.Range("D2").Value = "'245.35"
$ws.Range("D3").Value = "'22.01"
$ws.Range("D4").Value = "'5.341"
$ws.Range("D5").Value = "'0.05960"
$ws.Range("D6").Value = "'3.397"
$ws.Range("D7").Value = "'6.381"
$ws.Range("D8").Value = "'0.8084"
$ws.Range("D9").Value = "'0.9642"
$ws.Range("D11").Value = "'0.07409"
$ws.Range("D12").Value = "'0.03401"
$ws.Range("D13").Value = "'0.03072"
$ws.Range("D14").Value = "'0.09408"
$ws.Range("D15").Value = "'3.993"
$ws.Range("D16").Value = "'0.001594"
$ws.Range("D17").Value = "'0.04802"
$ws.Range("D18").Value = "'0.0005915"
$ws.Range("D19").Value = "'0.006208"
$ws.Range("D20").Value = "'0.005137"
$ws.Range("D21").Value = "'0.0009880"
$ws.Range("D23").Value = "'3.745"
$ws.Range("D27").Value = "'0.0002463"
$ws.Range("D44").Value = "'0.005323"
$ws.Range("D45").Value = "'0.00005314"
$ws.Range("D47").Value = "'0.8506"
$ws.Range("D48").Value = "'0.03738"

# --- Column E (Volume(1h)) text updates ---
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

# --- Rows 41-43 reshuffle (Coin / Link / Price / Volume(1h)) ---
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1072"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002692"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003035"
$ws.Range("E43").Value = "42KickTokenKICK"

# --- Column G (Hora) updates: 5 -> 6 for all data rows ---
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 7).Value = "'6"
}
